$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Il4"
$ws.Range("C2").Value = "Il13ra2"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1647103333333333
$ws.Range("H2").Value = 0.494131
$ws.Range("I2").Value = 0.05373200903458847
$ws.Range("J2").Value = 0.05373200903458848
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.900950666666667
$ws.Range("N2").Value = 5.702852
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.3131062179568889
$ws.Range("R2").Value = 2.817955961612
$ws.Range("S2").Value = 0.05373200903458847
$ws.Range("T2").Value = 0.05373200903458848

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Il4"
$ws.Range("C3").Value = "Il13ra2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.6783136666666666
$ws.Range("H3").Value = 2.034941
$ws.Range("I3").Value = 0.2212803248467603
$ws.Range("J3").Value = 0.2212803248467603
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.900950666666667
$ws.Range("N3").Value = 5.702852
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 1.289440816859111
$ws.Range("R3").Value = 11.604967351732
$ws.Range("S3").Value = 0.2212803248467603
$ws.Range("T3").Value = 0.2212803248467603

# Row 4
$ws.Range("A4").Value = "M1"
$ws.Range("B4").Value = "Il4"
$ws.Range("C4").Value = "Il13ra2"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.3320596666666667
$ws.Range("H4").Value = 0.996179
$ws.Range("I4").Value = 0.1083249159191941
$ws.Range("J4").Value = 0.1083249159191941
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.900950666666667
$ws.Range("N4").Value = 5.702852
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.6312290447231111
$ws.Range("R4").Value = 5.681061402508
$ws.Range("S4").Value = 0.1083249159191941
$ws.Range("T4").Value = 0.1083249159191941

# Row 5
$ws.Range("A5").Value = "M2"
$ws.Range("B5").Value = "Il4"
$ws.Range("C5").Value = "Il13ra2"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.272073
$ws.Range("H5").Value = 0.816219
$ws.Range("I5").Value = 0.08875599118898181
$ws.Range("J5").Value = 0.08875599118898181
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.900950666666667
$ws.Range("N5").Value = 5.702852
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 0.517197350732
$ws.Range("R5").Value = 4.654776156588
$ws.Range("S5").Value = 0.08875599118898181
$ws.Range("T5").Value = 0.08875599118898181

# Row 6
$ws.Range("A6").Value = "Neutro"
$ws.Range("B6").Value = "Il4"
$ws.Range("C6").Value = "Il13ra2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.2558096666666667
$ws.Range("H6").Value = 0.767429
$ws.Range("I6").Value = 0.08345054643688657
$ws.Range("J6").Value = 0.08345054643688658
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.900950666666667
$ws.Range("N6").Value = 5.702852
$ws.Range("O6").Value = 1
$ws.Range("P6").Value = 1
$ws.Range("Q6").Value = 0.4862815563897778
$ws.Range("R6").Value = 4.376534007508
$ws.Range("S6").Value = 0.08345054643688657
$ws.Range("T6").Value = 0.08345054643688658

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Il4"
$ws.Range("C7").Value = "Il13ra2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.362438
$ws.Range("H7").Value = 4.087314
$ws.Range("I7").Value = 0.4444562125735887
$ws.Range("J7").Value = 0.4444562125735887
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.900950666666667
$ws.Range("N7").Value = 5.702852
$ws.Range("O7").Value = 1
$ws.Range("P7").Value = 1
$ws.Range("Q7").Value = 2.589927424392
$ws.Range("R7").Value = 23.309346819528
$ws.Range("S7").Value = 0.4444562125735887
$ws.Range("T7").Value = 0.4444562125735887

